# Auto-generated edit script applying numeric corrections to the
# "Asura_Profits" workbook (currentAveragePrice / LevePrice / LeveProfit
# columns) per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3802
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 3802
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 4144.3125
$ws.Range("I76").Value = 4040.9
$ws.Range("K76").Value = 4040.9
$ws.Range("M76").Value = -3725.9
$ws.Range("H79").Value = 4144.3125
$ws.Range("I79").Value = 4040.9
$ws.Range("K79").Value = 4040.9
$ws.Range("M79").Value = -2948.9
$ws.Range("H132").Value = 1415.8948
$ws.Range("I132").Value = 1299.7747
$ws.Range("J132").Value = 3064.8
$ws.Range("K132").Value = 3899.3241
$ws.Range("L132").Value = 9194.400000000001
$ws.Range("M132").Value = -1369.3241
$ws.Range("N132").Value = -14254.4
$ws.Range("H137").Value = 1304.5
$ws.Range("I137").Value = 1111.102
$ws.Range("J137").Value = 3199.8
$ws.Range("K137").Value = 3333.306
$ws.Range("L137").Value = 9599.400000000001
$ws.Range("M137").Value = -783.3060000000005
$ws.Range("N137").Value = -14699.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15430.928
$ws.Range("I32").Value = 13774.75
$ws.Range("K32").Value = 13774.75
$ws.Range("M32").Value = -13487.75
$ws.Range("H43").Value = 500004540
$ws.Range("J43").Value = 500004540
$ws.Range("L43").Value = 500004540
$ws.Range("N43").Value = -500005166
$ws.Range("H63").Value = 13800
$ws.Range("I63").Value = 24000
$ws.Range("K63").Value = 24000
$ws.Range("M63").Value = -23314
$ws.Range("H66").Value = 13800
$ws.Range("I66").Value = 24000
$ws.Range("K66").Value = 120000
$ws.Range("M66").Value = -116568
$ws.Range("H74").Value = 1217.3695
$ws.Range("I74").Value = 962.2973
$ws.Range("J74").Value = 2266
$ws.Range("K74").Value = 962.2973
$ws.Range("L74").Value = 2266
$ws.Range("M74").Value = -88.29729999999995
$ws.Range("N74").Value = -4014
$ws.Range("H77").Value = 1217.3695
$ws.Range("I77").Value = 962.2973
$ws.Range("J77").Value = 2266
$ws.Range("K77").Value = 4811.4865
$ws.Range("L77").Value = 11330
$ws.Range("M77").Value = -443.4865
$ws.Range("N77").Value = -20066

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5088.5
$ws.Range("I105").Value = 4106.2
$ws.Range("K105").Value = 4106.2
$ws.Range("M105").Value = -2359.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3259.2156
$ws.Range("I31").Value = 1896.8462
$ws.Range("J31").Value = 4676.08
$ws.Range("K31").Value = 1896.8462
$ws.Range("L31").Value = 4676.08
$ws.Range("M31").Value = -1601.8462
$ws.Range("N31").Value = -5266.08
$ws.Range("H34").Value = 3259.2156
$ws.Range("I34").Value = 1896.8462
$ws.Range("J34").Value = 4676.08
$ws.Range("K34").Value = 1896.8462
$ws.Range("L34").Value = 4676.08
$ws.Range("M34").Value = -1694.8462
$ws.Range("N34").Value = -5080.08
$ws.Range("H58").Value = 1188.975
$ws.Range("I58").Value = 1204.6538
$ws.Range("J58").Value = 1159.8572
$ws.Range("K58").Value = 1204.6538
$ws.Range("L58").Value = 1159.8572
$ws.Range("M58").Value = -1001.6538
$ws.Range("N58").Value = -1565.8572
$ws.Range("H132").Value = 1959.8334
$ws.Range("I132").Value = 1625.8148
$ws.Range("K132").Value = 4877.4444
$ws.Range("M132").Value = -2347.4444
$ws.Range("H134").Value = 1018.2
$ws.Range("I134").Value = 962.2564
$ws.Range("J134").Value = 3200
$ws.Range("K134").Value = 2886.7692
$ws.Range("L134").Value = 9600
$ws.Range("M134").Value = -351.7691999999997
$ws.Range("N134").Value = -14670
$ws.Range("H136").Value = 1188.975
$ws.Range("I136").Value = 1204.6538
$ws.Range("J136").Value = 1159.8572
$ws.Range("K136").Value = 3613.9614
$ws.Range("L136").Value = 3479.5716
$ws.Range("M136").Value = -1063.9614
$ws.Range("N136").Value = -8579.571599999999
$ws.Range("H141").Value = 43195.6
$ws.Range("J141").Value = 43994.25
$ws.Range("L141").Value = 43994.25
$ws.Range("N141").Value = -54354.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 83334090
$ws.Range("I22").Value = 166666670
$ws.Range("J22").Value = 1498
$ws.Range("K22").Value = 500000010
$ws.Range("L22").Value = 4494
$ws.Range("M22").Value = -499999841
$ws.Range("N22").Value = -4832
$ws.Range("H26").Value = 167.5
$ws.Range("I26").Value = 148.57143
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 445.71429
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = -157.71429
$ws.Range("N26").Value = -1476
$ws.Range("H27").Value = 83334090
$ws.Range("I27").Value = 166666670
$ws.Range("J27").Value = 1498
$ws.Range("K27").Value = 500000010
$ws.Range("L27").Value = 4494
$ws.Range("M27").Value = -499999908
$ws.Range("N27").Value = -4698
$ws.Range("H34").Value = 811.9167
$ws.Range("I34").Value = 147.8
$ws.Range("J34").Value = 1286.2858
$ws.Range("K34").Value = 443.4
$ws.Range("L34").Value = 3858.8574
$ws.Range("M34").Value = -359.4
$ws.Range("N34").Value = -4026.8574
$ws.Range("H131").Value = 11366190
$ws.Range("J131").Value = 12659642
$ws.Range("L131").Value = 37978926
$ws.Range("N131").Value = -37989006
$ws.Range("H134").Value = 4526.857
$ws.Range("J134").Value = 7174.846
$ws.Range("L134").Value = 21524.538
$ws.Range("N134").Value = -31664.538

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8574857
$ws.Range("J3").Value = 36676336
$ws.Range("L3").Value = 36676336
$ws.Range("N3").Value = -36676568
$ws.Range("H70").Value = 5706.484
$ws.Range("I70").Value = 4993.6924
$ws.Range("J70").Value = 6221.278
$ws.Range("K70").Value = 4993.6924
$ws.Range("L70").Value = 6221.278
$ws.Range("M70").Value = -4723.6924
$ws.Range("N70").Value = -6761.278
$ws.Range("H73").Value = 5706.484
$ws.Range("I73").Value = 4993.6924
$ws.Range("J73").Value = 6221.278
$ws.Range("K73").Value = 4993.6924
$ws.Range("L73").Value = 6221.278
$ws.Range("M73").Value = -4057.6924
$ws.Range("N73").Value = -8093.278
$ws.Range("H80").Value = 3200.5
$ws.Range("J80").Value = 3625
$ws.Range("L80").Value = 3625
$ws.Range("N80").Value = -5621
$ws.Range("H83").Value = 3200.5
$ws.Range("J83").Value = 3625
$ws.Range("L83").Value = 18125
$ws.Range("N83").Value = -28109
$ws.Range("H103").Value = 80600
$ws.Range("J103").Value = 80600
$ws.Range("L103").Value = 80600
$ws.Range("N103").Value = -82944
$ws.Range("H132").Value = 1526.4324
$ws.Range("I132").Value = 987.9583
$ws.Range("J132").Value = 2520.5386
$ws.Range("K132").Value = 2963.8749
$ws.Range("L132").Value = 7561.6158
$ws.Range("M132").Value = -433.8748999999998
$ws.Range("N132").Value = -12621.6158

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 122527.164
$ws.Range("J69").Value = 122527.164
$ws.Range("L69").Value = 122527.164
$ws.Range("N69").Value = -124149.164
$ws.Range("H72").Value = 122527.164
$ws.Range("J72").Value = 122527.164
$ws.Range("L72").Value = 367581.492
$ws.Range("N72").Value = -375693.492
$ws.Range("H136").Value = 3743.923
$ws.Range("I136").Value = 3701.681
$ws.Range("J136").Value = 4141
$ws.Range("K136").Value = 11105.043
$ws.Range("L136").Value = 12423
$ws.Range("M136").Value = -8555.043
$ws.Range("N136").Value = -17523

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1264.6714
$ws.Range("I132").Value = 1144.4386
$ws.Range("J132").Value = 1791.8462
$ws.Range("K132").Value = 3433.3158
$ws.Range("L132").Value = 5375.5386
$ws.Range("M132").Value = -903.3157999999999
$ws.Range("N132").Value = -10435.5386
$ws.Range("H136").Value = 2106.8708
$ws.Range("I136").Value = 2061.3572
$ws.Range("J136").Value = 2531.6667
$ws.Range("K136").Value = 6184.071599999999
$ws.Range("L136").Value = 7595.000100000001
$ws.Range("M136").Value = -3634.071599999999
$ws.Range("N136").Value = -12695.0001

